$d = $word.ActiveDocument

$replacements = @(
    @{old = "(Ref-J9DH4K)"; new = "(Pearse et al. 117)"},
    @{old = "(Ref-29F0B8)"; new = "(Ref-s918663)"},
    @{old = "(Ref-J7X8K2)"; new = "(Johnson)"},
    @{old = "(Ref-J6DH3Y)"; new = "(Ref-s096532)"},
    @{old = "(Ref-AB12CD)"; new = "(Smith)"},
    @{old = "(Ref-EF34GH)"; new = "(Smith)"},
    @{old = "(Ref-J7X2B9)"; new = "(Ref-u321387)"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
